$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column (shifts existing Project Name..Officer data from A:N to B:O)
$ws.Columns.Item(1).Insert()

# New "Project ID" column
$ws.Range("A1").Value = "Project ID"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2

# New project row: Melville Park / Tampines
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Melville Park"
$ws.Range("C4").Value = "Tampines"
$ws.Range("D4").Value = "2-ROOM"
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 1000000
$ws.Range("G4").Value = "3-ROOM"
$ws.Range("H4").Value = 10
$ws.Range("I4").Value = 2000000
$ws.Range("J4").Value = 45775
$ws.Range("K4").Value = 45805
$ws.Range("L4").Value = "T8765432F"
$ws.Range("M4").Value = 10
$ws.Range("O4").Value = "Hidden"

$ws.Range("L4").Select()
